# corrected data cleaning for pre/post/total fixation data
#
# The exported dataframe used to contain a spurious duplicated header row
# (old row 2) and a header row that still carried the pandas index label
# "Unnamed: 0" plus a bold/centered/bordered header style. The corrected
# export drops that duplicate row, clears the header styling/label, and
# ends with two trailing blank rows (rows 9-10) instead of the data
# running out at row 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicated header row (old row 2) - this shifts every
# row below it up by one (old rows 3..9 become new rows 2..8).
$ws.Rows.Item(2).Delete()

# The header row (row 1) no longer has the bold/centered/bordered style,
# and A1 no longer holds the "Unnamed: 0" pandas index label.
$ws.Range("A1:W1").ClearFormats()
$ws.Range("A1").Value = ""

# The cleaned export now keeps two trailing blank rows (9 and 10) after
# the last data row ("First fixation duration (ms)", now row 8). Touch
# their formatting (a no-op) so the used range/dimension grows to
# include them as blank rows rather than stopping at row 8.
$ws.Range("A9:W10").Font.Bold = $false

$wb.Save()
